# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 10:03"

# Row 37 - Polonia
$ws.Range("D37").Value = 80302
$ws.Range("E37").Value = 38364

# Row 63 - Armenia
$ws.Range("B63").Value = 56451
$ws.Range("C63").Value = 715
$ws.Range("D63").Value = 45824
$ws.Range("E63").Value = 9607
$ws.Range("G63").Value = 4
$ws.Range("H63").Value = 1020

# Row 75 - Afganistan
$ws.Range("B75").Value = 39799
$ws.Range("C75").Value = 96
$ws.Range("D75").Value = 33114
$ws.Range("E75").Value = 5208
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 1477

# Row 76 - Hungria
$ws.Range("B76").Value = 37664
$ws.Range("C76").Value = 1068
$ws.Range("D76").Value = 10848
$ws.Range("E76").Value = 25862
$ws.Range("G76").Value = 21
$ws.Range("H76").Value = 954

# Row 103 - Georgia
$ws.Range("B103").Value = 11794
$ws.Range("C103").Value = 523
$ws.Range("D103").Value = 6327
$ws.Range("E103").Value = 5382
$ws.Range("G103").Value = 7
$ws.Range("H103").Value = 85

# Rows 140/141 - Estonia moves above Somalia in the ranking, with refreshed
# data for Estonia while Somalia keeps its previous figures.
$ws.Range("A140").Value = "Estonia"
$ws.Range("B140").Value = 3865
$ws.Range("C140").Value = 19
$ws.Range("D140").Value = 2958
$ws.Range("E140").Value = 839
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 68

$ws.Range("A141").Value = "Somalia"
$ws.Range("B141").Value = 3847
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 3079
$ws.Range("E141").Value = 669
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 99

# Row 152 - Letonia
$ws.Range("B152").Value = 2670
$ws.Range("C152").Value = 74
$ws.Range("E152").Value = 1308
